$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update runs (C), balls (D), fours (E) columns for rows 2-6
# to reflect the "updated activity till excel form" values.
# Force text format so values stay stored as text (matching the
# original cells, which were all text/string typed).

$textRange = $ws.Range("C2:E6")
$textRange.NumberFormat = "@"

$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "2"
$ws.Range("E2").Value = "0"

$ws.Range("C3").Value = "23"
$ws.Range("D3").Value = "23"
$ws.Range("E3").Value = "2"

$ws.Range("C4").Value = "6"
$ws.Range("D4").Value = "4"
$ws.Range("E4").Value = "1"

$ws.Range("C5").Value = "5"
$ws.Range("D5").Value = "7"
$ws.Range("E5").Value = "0"

$ws.Range("C6").Value = "2"
$ws.Range("D6").Value = "3"
